# Update probability-matrix cells on Sheet1 with refreshed values
# (games pulled March 7) per the team-specific transition matrix.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1616541353383459
$ws.Range("C2").Value = 0.6390977443609023
$ws.Range("J2").Value = 0.0112781954887218
$ws.Range("P2").Value = 0.1428571428571428
$ws.Range("S2").Value = 0.04511278195488722
$ws.Range("B3").Value = 0.005494505494505495
$ws.Range("C3").Value = 0.04945054945054945
$ws.Range("J3").Value = 0.005494505494505495
$ws.Range("P3").Value = 0.8076923076923077
$ws.Range("S3").Value = 0.1318681318681319
$ws.Range("J4").Value = 0.02439024390243903
$ws.Range("P4").Value = 0.7317073170731707
$ws.Range("S4").Value = 0.2439024390243902
$ws.Range("B6").Value = 0.04219409282700422
$ws.Range("D6").Value = 0.004219409282700422
$ws.Range("F6").Value = 0.08016877637130802
$ws.Range("J6").Value = 0.2658227848101266
$ws.Range("O6").Value = 0.02109704641350211
$ws.Range("Q6").Value = 0.1687763713080169
$ws.Range("R6").Value = 0.06329113924050633
$ws.Range("S6").Value = 0.3544303797468354
$ws.Range("B7").Value = 0.08
$ws.Range("D7").Value = 0.02285714285714286
$ws.Range("F7").Value = 0.06857142857142857
$ws.Range("J7").Value = 0.12
$ws.Range("O7").Value = 0.01142857142857143
$ws.Range("Q7").Value = 0.2114285714285714
$ws.Range("R7").Value = 0.1142857142857143
$ws.Range("S7").Value = 0.3714285714285714
$ws.Range("B8").Value = 0.1176470588235294
$ws.Range("D8").Value = 0.01890756302521008
$ws.Range("E8").Value = 0.002100840336134454
$ws.Range("F8").Value = 0.06092436974789916
$ws.Range("J8").Value = 0.09033613445378151
$ws.Range("O8").Value = 0.01050420168067227
$ws.Range("Q8").Value = 0.1554621848739496
$ws.Range("R8").Value = 0.1113445378151261
$ws.Range("S8").Value = 0.4327731092436975
$ws.Range("B9").Value = 0.07303370786516854
$ws.Range("D9").Value = 0.02808988764044944
$ws.Range("E9").Value = 0.005617977528089887
$ws.Range("F9").Value = 0.07303370786516854
$ws.Range("J9").Value = 0.101123595505618
$ws.Range("O9").Value = 0.0449438202247191
$ws.Range("Q9").Value = 0.2247191011235955
$ws.Range("R9").Value = 0.09550561797752809
$ws.Range("S9").Value = 0.3539325842696629
$ws.Range("B10").Value = 0.1098712446351931
$ws.Range("D10").Value = 0.01888412017167382
$ws.Range("F10").Value = 0.0703862660944206
$ws.Range("J10").Value = 0.09012875536480687
$ws.Range("O10").Value = 0.01545064377682403
$ws.Range("Q10").Value = 0.2206008583690987
$ws.Range("R10").Value = 0.1072961373390558
$ws.Range("S10").Value = 0.367381974248927
$ws.Range("G11").Value = 0.144
$ws.Range("J11").Value = 0.07199999999999999
$ws.Range("K11").Value = 0.18
$ws.Range("L11").Value = 0.58
$ws.Range("S11").Value = 0.024
$ws.Range("G12").Value = 0.756578947368421
$ws.Range("J12").Value = 0.1644736842105263
$ws.Range("L12").Value = 0.05263157894736842
$ws.Range("S12").Value = 0.02631578947368421
$ws.Range("G13").Value = 0.5714285714285714
$ws.Range("J13").Value = 0.3673469387755102
$ws.Range("S13").Value = 0.06122448979591837
$ws.Range("F15").Value = 0.02100840336134454
$ws.Range("H15").Value = 0.180672268907563
$ws.Range("I15").Value = 0.06302521008403361
$ws.Range("J15").Value = 0.3445378151260504
$ws.Range("K15").Value = 0.05882352941176471
$ws.Range("M15").Value = 0.01260504201680672
$ws.Range("O15").Value = 0.1050420168067227
$ws.Range("S15").Value = 0.2142857142857143
$ws.Range("F16").Value = 0.05213270142180094
$ws.Range("H16").Value = 0.1895734597156398
$ws.Range("I16").Value = 0.06161137440758294
$ws.Range("J16").Value = 0.4075829383886256
$ws.Range("K16").Value = 0.08530805687203792
$ws.Range("M16").Value = 0.01895734597156398
$ws.Range("N16").Value = 0.004739336492890996
$ws.Range("O16").Value = 0.07109004739336493
$ws.Range("S16").Value = 0.1090047393364929
$ws.Range("F17").Value = 0.02643171806167401
$ws.Range("H17").Value = 0.1982378854625551
$ws.Range("I17").Value = 0.1013215859030837
$ws.Range("J17").Value = 0.4140969162995595
$ws.Range("K17").Value = 0.07709251101321586
$ws.Range("M17").Value = 0.01541850220264317
$ws.Range("O17").Value = 0.05506607929515418
$ws.Range("S17").Value = 0.1123348017621145
$ws.Range("F18").Value = 0.04291845493562232
$ws.Range("H18").Value = 0.1459227467811159
$ws.Range("I18").Value = 0.09871244635193133
$ws.Range("J18").Value = 0.4506437768240343
$ws.Range("K18").Value = 0.08583690987124463
$ws.Range("M18").Value = 0.01716738197424893
$ws.Range("O18").Value = 0.0815450643776824
$ws.Range("S18").Value = 0.07725321888412018
$ws.Range("F19").Value = 0.02144082332761578
$ws.Range("H19").Value = 0.2332761578044597
$ws.Range("I19").Value = 0.07204116638078903
$ws.Range("J19").Value = 0.3447684391080618
$ws.Range("K19").Value = 0.09862778730703259
$ws.Range("M19").Value = 0.02658662092624357
$ws.Range("N19").Value = 0.0008576329331046312
$ws.Range("O19").Value = 0.07375643224699828
$ws.Range("S19").Value = 0.1286449399656947
